$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "LED Vario 160"
$ws.Range("H1").Value = "PT Kytaco Japan"

$ws.Range("F5").Select()
